$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("A6").Value = "2024-08-04 22:11:04"
$ws.Range("C6").Value = 12
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 0
$ws.Range("N6").Value = 10
$ws.Range("O6").Value = 6
$ws.Range("P6").Value = 3
$ws.Range("R6").Value = 5
$ws.Range("T6").Value = 20
$ws.Range("U6").Value = 0.5
$ws.Range("V6").Value = "D:\Repositorio\jonatha1992\Predictor_ruleta\Data\Electromecanica.xlsx"
$ws.Range("X6").Value = "No es Simulación"
$ws.Range("Y6").Value = 24

# Row 7
$ws.Range("A7").Value = "2024-08-04 22:13:50"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("N7").Value = 10
$ws.Range("O7").Value = 6
$ws.Range("P7").Value = 3
$ws.Range("R7").Value = 5
$ws.Range("T7").Value = 20
$ws.Range("U7").Value = 0
$ws.Range("V7").Value = "D:\Repositorio\jonatha1992\Predictor_ruleta\Data\Crupier.xlsx"
$ws.Range("X7").Value = "No es Simulación"
$ws.Range("Y7").Value = 0
